$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.25"
$ws.Range("E2").Value = "'0.33%"

$ws.Range("D3").Value = "'44.27"
$ws.Range("E3").Value = "'0.76%"

$ws.Range("D4").Value = "'5.573"
$ws.Range("E4").Value = "'1.79%"

$ws.Range("E5").Value = "'0.20%"

$ws.Range("D6").Value = "'1.977"
$ws.Range("E6").Value = "'4.75%"

$ws.Range("D7").Value = "'4.326"
$ws.Range("E7").Value = "'0.82%"

$ws.Range("D8").Value = "'0.9523"
$ws.Range("E8").Value = "'1.59%"

$ws.Range("E9").Value = "'-4.76%"

$ws.Range("D10").Value = "'0.1169"
$ws.Range("E10").Value = "'-1.88%"

$ws.Range("D11").Value = "'0.1859"
$ws.Range("E11").Value = "'-1.86%"

$ws.Range("D12").Value = "'10.30"
$ws.Range("E12").Value = "'19.61%"

$ws.Range("D13").Value = "'0.09849"
$ws.Range("E13").Value = "'2.42%"

$ws.Range("D14").Value = "'0.04717"
$ws.Range("E14").Value = "'15.28%"

$ws.Range("D15").Value = "'0.1067"
$ws.Range("E15").Value = "'-0.06%"

$ws.Range("D16").Value = "'0.001291"
$ws.Range("E16").Value = "'0.94%"

$ws.Range("D17").Value = "'0.04228"
$ws.Range("E17").Value = "'-3.08%"

$ws.Range("D18").Value = "'0.005858"
$ws.Range("E18").Value = "'-1.84%"

$ws.Range("D19").Value = "'3.376"
$ws.Range("E19").Value = "'-5.39%"

$ws.Range("D20").Value = "'0.3473"
$ws.Range("E20").Value = "'-0.73%"

$ws.Range("D21").Value = "'0.1409"
$ws.Range("E21").Value = "'4.54%"

$ws.Range("D22").Value = "'0.2507"
$ws.Range("E22").Value = "'0.61%"

$ws.Range("D23").Value = "'0.001258"
$ws.Range("E23").Value = "'1.91%"

$ws.Range("D24").Value = "'0.004330"
$ws.Range("E24").Value = "'0.60%"

$ws.Range("D25").Value = "'0.0001191"
$ws.Range("E25").Value = "'-3.43%"

$ws.Range("E26").Value = "'-0.51%"

$ws.Range("D38").Value = "'0.02663"
$ws.Range("E38").Value = "'-0.12%"

$ws.Range("D39").Value = "'0.05539"
$ws.Range("E39").Value = "'1.78%"

$ws.Range("D40").Value = "'0.007575"
$ws.Range("E40").Value = "'-1.37%"

$ws.Range("D41").Value = "'0.1408"
$ws.Range("E41").Value = "'1.33%"

$ws.Range("D42").Value = "'0.008081"
$ws.Range("E42").Value = "'-17.13%"

$ws.Range("D43").Value = "'0.002017"
$ws.Range("E43").Value = "'-3.52%"

$ws.Range("D44").Value = "'0.008902"
$ws.Range("E44").Value = "'-10.11%"

$ws.Range("D45").Value = "'0.00007272"
$ws.Range("E45").Value = "'2.19%"

$ws.Range("E46").Value = "'-0.20%"

$ws.Range("D47").Value = "'0.004765"
$ws.Range("E47").Value = "'33.79%"

$ws.Range("D48").Value = "'0.002272"
$ws.Range("E48").Value = "'-0.18%"

$ws.Range("D49").Value = "'0.00002102"
$ws.Range("E49").Value = "'-0.20%"

$ws.Range("D50").Value = "'0.0002002"
$ws.Range("E50").Value = "'-0.20%"
